$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B2 value (BRUNO370 -> BRUNO383)
$ws.Range("B2").Value = "BRUNO383"

# Enter the new rows' data in the same order the shared-string table was
# built (reconstructed from the target XML) so the resulting sharedStrings
# ordering matches the target exactly.
$ws.Range("A3").Value = "2º"
$ws.Range("A4").Value = "3º"
$ws.Range("E3").Value = "Marcela"
$ws.Range("F3").Value = "Matos"
$ws.Range("C3").Value = "marcelamatos@teste.com"
$ws.Range("B3").Value = "Marcela74"
$ws.Range("B4").Value = "BrUn-^:1SW"
$ws.Range("C4").Value = "blumenal@teste.com"
$ws.Range("D3").Value = "Marcela.1234"
$ws.Range("D4").Value = "paçoca123A"
$ws.Range("E4").Value = "Blumenal"
$ws.Range("F4").Value = "Souza"
$ws.Range("G3").Value = "55 11 8542 3671"
$ws.Range("G4").Value = "57 21 8545 3535"
$ws.Range("H3").Value = "Bahamas"
$ws.Range("H4").Value = "French Polynesia"
$ws.Range("I3").Value = "Galapz"
$ws.Range("I4").Value = "Martito"
$ws.Range("J3").Value = "Futton Six"
$ws.Range("J4").Value = "Street Max"
$ws.Range("K3").Value = "AS"
$ws.Range("K4").Value = "LS"
$ws.Range("L3").Value = "00 555 111 3"
$ws.Range("L4").Value = "40852 41"

# Add the mail hyperlinks on the e-mail column for the two new rows before
# re-applying the row-2 formatting, so the hyperlink cells end up sharing
# the very same cell style (border/font/alignment) as C2 instead of the
# default "fresh hyperlink" style.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:marcelamatos@teste.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:blumenal@teste.com")

# Copy row 2's formatting (borders/fonts/alignment/row height) down onto
# the two new rows.
$ws.Range("A2:L2").Copy()
$ws.Range("A3:L4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows("3:4").RowHeight = 15.75

# Match the final selection/active cell recorded in the sheet view.
$ws.Range("L4").Select()
